$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.133.42'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '2.308.23'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '301.64'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.46%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '100.10'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.85%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.513'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +1.16%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.512'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.08%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.44'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +8.15%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0793'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('E12').Value = '  +0.62%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '17.76'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +3.78%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.92'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.16%  '
$ws.Range('D15').Value = '2.663.70'
$ws.Range('E15').Value = '  +0.36%  '
$ws.Range('D16').Value = '2.269.57'
$ws.Range('E16').Value = '  -1.39%  '
$ws.Range('E17').Value = '  -0.98%  '
$ws.Range('D18').Value = '43.028.87'
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('E19').Value = '  +10.00%  '
$ws.Range('D20').Value = '0.0₃0905'
$ws.Range('E20').Value = '  +0.45%  '
$ws.Range('E21').Value = '  +1.22%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '67.99'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.62%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '235.80'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.20'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +7.72%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.46'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '25.14'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.94%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '169.93'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.88%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '34.46'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.81%  '
$ws.Range('E31').Value = '  +0.54%  '
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('E33').Value = '  +2.19%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '17.73'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +6.05%  '
$ws.Range('E35').Value = '  +0.86%  '
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('E37').Value = '  -0.70%  '
$ws.Range('E38').Value = '  +1.46%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.79'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.37%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.82'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.62%  '
$ws.Range('E41').Value = '  +0.64%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0292'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +3.08%  '
$ws.Range('D43').Value = '1.986.66'
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.26'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -4.80%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '10.20'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.84%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '17.63'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.90'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.81%  '
$ws.Range('E48').Value = '  +4.39%  '
$ws.Range('D50').Value = '2.531.50'
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '70.81'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.39%  '
